$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title rename: "Initial experiments" -> "Experiments" ---
$ws.Range("A1").Value = "Experiments"

# --- New header cells for the first table (Precision, Recall, Distance) ---
$ws.Range("M2").Value = "Precision"
$ws.Range("N2").Value = "Recall"
$ws.Range("O2").Value = "Distance"

# --- New data row 9 (7th run) in the first table ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 200
$ws.Range("C9").Value = 200
$ws.Range("D9").Value = "rgb"
$ws.Range("E9").Value = 1600
$ws.Range("F9").Value = 1500
$ws.Range("G9").Value = 100
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 10
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 0.61
$ws.Range("L9").Value = 0.65
$ws.Range("L9").NumberFormat = "0%"
$ws.Range("M9").Value = 0.24
$ws.Range("N9").Value = 0.38
$ws.Range("O9").Value = "Manhattan"

# --- New data row 15 (2nd run) in the second (kmeans distance) table ---
$ws.Range("A15").Value = 2
$ws.Range("B15").Value = 200
$ws.Range("C15").Value = 200
$ws.Range("D15").Value = "rgb"
$ws.Range("E15").Value = 1600
$ws.Range("F15").Value = 1500
$ws.Range("G15").Value = 100
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 10
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 0.61
$ws.Range("L15").NumberFormat = "0%"

# --- Four new comparison rows (22-25) for Setup 2, mirroring rows 18-21 ---
$ws.Range("A22").Value = 2
$ws.Range("B22").Value = 0.45
$ws.Range("B22").NumberFormat = "0%"
$ws.Range("C22").Value = "Euclidean"
$ws.Range("D22").Value = "sklearn"

$ws.Range("A23").Value = 2
$ws.Range("B23").Value = 0.64
$ws.Range("B23").NumberFormat = "0%"
$ws.Range("C23").Value = "Cosine"
$ws.Range("D23").Value = "nltk"

$ws.Range("A24").Value = 2
$ws.Range("B24").Value = 0.45
$ws.Range("B24").NumberFormat = "0%"
$ws.Range("C24").Value = "Euclidean"
$ws.Range("D24").Value = "nltk"

$ws.Range("A25").Value = 2
$ws.Range("B25").Value = 0.65
$ws.Range("B25").NumberFormat = "0%"
$ws.Range("C25").Value = "Manhattan"
$ws.Range("D25").Value = "nltk"

# --- Update view state: scroll so row 7 is the top visible row, then select E22 ---
$ws.Activate() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 7
} catch {
    # Scroll-position persistence may be unsupported by the host; selection below still applies.
}
$ws.Range("E22").Select() | Out-Null
